$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-27 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-28 Thursday", 2) | Out-Null
$d.Content.Find.Execute("521÷3=173, 2", $true, $false, $false, $false, $false, $true, 1, $false, "410÷2=205, 0", 2) | Out-Null
$d.Content.Find.Execute("669÷6=111, 3", $true, $false, $false, $false, $false, $true, 1, $false, "576÷2=288, 0", 2) | Out-Null
$d.Content.Find.Execute("708÷7=101, 1", $true, $false, $false, $false, $false, $true, 1, $false, "735÷3=245, 0", 2) | Out-Null
$d.Content.Find.Execute("967÷4=241, 3", $true, $false, $false, $false, $false, $true, 1, $false, "344÷2=172, 0", 2) | Out-Null
$d.Content.Find.Execute("255÷8=31, 7", $true, $false, $false, $false, $false, $true, 1, $false, "390÷9=43, 3", 2) | Out-Null
$d.Content.Find.Execute("788÷9=87, 5", $true, $false, $false, $false, $false, $true, 1, $false, "833÷5=166, 3", 2) | Out-Null
$d.Content.Find.Execute("853÷3=284, 1", $true, $false, $false, $false, $false, $true, 1, $false, "581÷9=64, 5", 2) | Out-Null
$d.Content.Find.Execute("110÷2=55, 0", $true, $false, $false, $false, $false, $true, 1, $false, "113÷7=16, 1", 2) | Out-Null
$d.Content.Find.Execute("465÷9=51, 6", $true, $false, $false, $false, $false, $true, 1, $false, "668÷3=222, 2", 2) | Out-Null
$d.Content.Find.Execute("573÷7=81, 6", $true, $false, $false, $false, $false, $true, 1, $false, "445÷6=74, 1", 2) | Out-Null
$d.Content.Find.Execute("755÷8=94, 3", $true, $false, $false, $false, $false, $true, 1, $false, "779÷9=86, 5", 2) | Out-Null
$d.Content.Find.Execute("117÷3=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "924÷8=115, 4", 2) | Out-Null
$d.Content.Find.Execute("581÷2=290, 1", $true, $false, $false, $false, $false, $true, 1, $false, "492÷7=70, 2", 2) | Out-Null
$d.Content.Find.Execute("198÷5=39, 3", $true, $false, $false, $false, $false, $true, 1, $false, "976÷9=108, 4", 2) | Out-Null
$d.Content.Find.Execute("101÷6=16, 5", $true, $false, $false, $false, $false, $true, 1, $false, "867÷2=433, 1", 2) | Out-Null
$d.Content.Find.Execute("497÷6=82, 5", $true, $false, $false, $false, $false, $true, 1, $false, "668÷3=222, 2", 2) | Out-Null
$d.Content.Find.Execute("494÷8=61, 6", $true, $false, $false, $false, $false, $true, 1, $false, "934÷4=233, 2", 2) | Out-Null
$d.Content.Find.Execute("167÷4=41, 3", $true, $false, $false, $false, $false, $true, 1, $false, "313÷7=44, 5", 2) | Out-Null
$d.Content.Find.Execute("305÷7=43, 4", $true, $false, $false, $false, $false, $true, 1, $false, "243÷8=30, 3", 2) | Out-Null
$d.Content.Find.Execute("997÷6=166, 1", $true, $false, $false, $false, $false, $true, 1, $false, "627÷2=313, 1", 2) | Out-Null
$d.Content.Find.Execute("645÷8=80, 5", $true, $false, $false, $false, $false, $true, 1, $false, "306÷4=76, 2", 2) | Out-Null
$d.Content.Find.Execute("214÷5=42, 4", $true, $false, $false, $false, $false, $true, 1, $false, "746÷9=82, 8", 2) | Out-Null
$d.Content.Find.Execute("249÷8=31, 1", $true, $false, $false, $false, $false, $true, 1, $false, "202÷2=101, 0", 2) | Out-Null
$d.Content.Find.Execute("194÷2=97, 0", $true, $false, $false, $false, $false, $true, 1, $false, "209÷5=41, 4", 2) | Out-Null
$d.Content.Find.Execute("190÷6=31, 4", $true, $false, $false, $false, $false, $true, 1, $false, "951÷3=317, 0", 2) | Out-Null
